# Generate Report for Handoff
# Adds a new source file (c8412e32-ecc8-4e5b-8951-35c821e836b9.md) to the
# localization-status report: one summary row on "Overview" and one detail
# row on each of the "zh-cn" / "de-de" sheets, pushing the existing
# ".localization-config" row down by one.

$wb = $excel.ActiveWorkbook

$newFile  = "c8412e32-ecc8-4e5b-8951-35c821e836b9.md"
$newZh    = "c8412e32-ecc8-4e5b-8951-35c821e836b9.6f9f00140ef3ccebc7c3a311aa08ad19b92a4aa1.zh-cn.xlf"
$newZhDt  = "2016-03-10 00:39:57"
$newDe    = "c8412e32-ecc8-4e5b-8951-35c821e836b9.6f9f00140ef3ccebc7c3a311aa08ad19b92a4aa1.de-de.xlf"
$newDeDt  = "2016-03-10 00:40:04"

$mdTarget     = "https://github.com/OpenLocalizationTest/oltest/blob/354b383d189e76cd012c8dabb5f3b3c855ed5a56/e2e/c8412e32-ecc8-4e5b-8951-35c821e836b9.md"
$cfgTarget    = "https://github.com/OpenLocalizationTest/oltest/blob/354b383d189e76cd012c8dabb5f3b3c855ed5a56/.localization-config"
$zhXlfTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3cf5a2d2e467bfc3bd5ab58fdab442aac6cf3621/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZh"
$deXlfTarget  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc50117e718b485c76d08edb409140b4604778b7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDe"
$mdTargetRepeat = "https://github.com/OpenLocalizationTest/oltest/blob/354b383d189e76cd012c8dabb5f3b3c855ed5a56/e2e/8bb8803f-22c8-4a89-8ba5-07e846cd53e7.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Duplicate row 2 (values + styles) into a fresh row 3, pushing the old
# row 3 (".localization-config") down to row 4.
$ws1.Range("A2:C2").Copy()
$ws1.Range("A3").Insert()

# Overwrite the new row 3 with the new file's data.
$ws1.Range("A3").Value = $newFile
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

# Rebuild the hyperlinks for this sheet (old ones don't track the row shift).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdTargetRepeat, "", "", "8bb8803f-22c8-4a89-8ba5-07e846cd53e7.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdTarget, "", "", $newFile)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $cfgTarget, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2:I2").Copy()
$ws2.Range("A3").Insert()

$ws2.Range("A3").Value = $newFile
$ws2.Range("C3").Value = $newZh
$ws2.Range("D3").Value = $newZhDt

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdTargetRepeat, "", "", "8bb8803f-22c8-4a89-8ba5-07e846cd53e7.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3cf5a2d2e467bfc3bd5ab58fdab442aac6cf3621/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8bb8803f-22c8-4a89-8ba5-07e846cd53e7.1379be1f6173d0999b12c3eb73d1a151b1ef53c6.zh-cn.xlf", "", "", "8bb8803f-22c8-4a89-8ba5-07e846cd53e7.1379be1f6173d0999b12c3eb73d1a151b1ef53c6.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdTarget, "", "", $newFile)
$ws2.Hyperlinks.Add($ws2.Range("C3"), $zhXlfTarget, "", "", $newZh)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $cfgTarget, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2:I2").Copy()
$ws3.Range("A3").Insert()

$ws3.Range("A3").Value = $newFile
$ws3.Range("C3").Value = $newDe
$ws3.Range("D3").Value = $newDeDt

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdTargetRepeat, "", "", "8bb8803f-22c8-4a89-8ba5-07e846cd53e7.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc50117e718b485c76d08edb409140b4604778b7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8bb8803f-22c8-4a89-8ba5-07e846cd53e7.1379be1f6173d0999b12c3eb73d1a151b1ef53c6.de-de.xlf", "", "", "8bb8803f-22c8-4a89-8ba5-07e846cd53e7.1379be1f6173d0999b12c3eb73d1a151b1ef53c6.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdTarget, "", "", $newFile)
$ws3.Hyperlinks.Add($ws3.Range("C3"), $deXlfTarget, "", "", $newDe)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $cfgTarget, "", "", ".localization-config")
